$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in the Price column (D) whose new value would otherwise be
# auto-parsed as a number (losing formatting like trailing zeros, e.g.
# "1.00" -> 1). Force those specific cells to Text first so the exact
# textual representation from the source data is preserved; cells that
# are naturally non-numeric (e.g. '67.988.35') don't need this.
$textForceRows = @(
    4, 5, 6, 10, 11, 12, 13, 14, 19, 20, 21, 22, 23, 24, 26, 27, 28, 30, 31, 32, 33, 35, 38, 39, 40, 41, 42, 43, 44, 45, 48, 49, 50, 51
)
foreach ($row in $textForceRows) {
    $ws.Cells.Item($row, 4).NumberFormat = "@"
}

$data = @{
    2 = @{ D = '67.988.35'; E = '  -2.80%  ' }
    3 = @{ D = '3.814.92'; E = '  +1.72%  ' }
    4 = @{ D = '1.00'; E = '  -0.18%  ' }
    5 = @{ D = '596.79'; E = '  -3.62%  ' }
    6 = @{ D = '174.43'; E = '  -3.79%  ' }
    7 = @{ D = '3.811.50'; E = '  +1.76%  ' }
    8 = @{ E = '  +0.04%  ' }
    9 = @{ E = '  -1.07%  ' }
    10 = @{ D = '0.160'; E = '  -3.80%  ' }
    11 = @{ D = '6.28'; E = '  -3.97%  ' }
    12 = @{ D = '0.464'; E = '  -3.82%  ' }
    13 = @{ D = '38.04'; E = '  -4.85%  ' }
    14 = @{ D = '0.0000245'; E = '  -4.33%  ' }
    15 = @{ D = '4.453.93'; E = '  +1.69%  ' }
    16 = @{ D = '3.815.46'; E = '  +1.45%  ' }
    17 = @{ D = '68.123.27'; E = '  -2.77%  ' }
    18 = @{ E = '  -4.47%  ' }
    19 = @{ D = '7.15'; E = '  -5.32%  ' }
    20 = @{ D = '16.29'; E = '  -2.08%  ' }
    21 = @{ D = '490.06'; E = '  -2.90%  ' }
    22 = @{ D = '9.22'; E = '  +0.15%  ' }
    23 = @{ D = '0.732'; E = '  +0.85%  ' }
    24 = @{ D = '84.72'; E = '  -2.27%  ' }
    25 = @{ E = '  -8.70%  ' }
    26 = @{ D = '0.0000138'; E = '  +2.82%  ' }
    27 = @{ D = '12.31'; E = '  -5.42%  ' }
    28 = @{ D = '10.23'; E = '  -9.47%  ' }
    29 = @{ E = '  +0.20%  ' }
    30 = @{ D = '2.92'; E = '  -0.54%  ' }
    31 = @{ D = '2.43'; E = '  -2.07%  ' }
    32 = @{ D = '32.75'; E = '  +6.69%  ' }
    33 = @{ D = '7.73'; E = '  -2.81%  ' }
    34 = @{ E = '  -4.31%  ' }
    35 = @{ D = '1.00'; E = '  -0.28%  ' }
    36 = @{ E = '  -4.31%  ' }
    37 = @{ E = '  -1.45%  ' }
    38 = @{ D = '5.78'; E = '  -6.28%  ' }
    39 = @{ D = '0.326'; E = '  -6.93%  ' }
    40 = @{ D = '449.36'; E = '  +2.15%  ' }
    41 = @{ D = '48.99'; E = '  -1.31%  ' }
    42 = @{ D = '1.99'; E = '  -3.77%  ' }
    43 = @{ D = '2.88'; E = '  -6.50%  ' }
    44 = @{ D = '8.28'; E = '  -4.03%  ' }
    45 = @{ D = '41.52'; E = '  -9.03%  ' }
    46 = @{ D = '2.829.80'; E = '  -5.12%  ' }
    48 = @{ D = '138.55'; E = '  -0.06%  ' }
    49 = @{ D = '0.0351' }
    50 = @{ D = '26.17'; E = '  -4.85%  ' }
    51 = @{ D = '23.11'; E = '  +6.49%  ' }
}

foreach ($row in $data.Keys) {
    $item = $data[$row]
    if ($item.ContainsKey("D")) {
        $ws.Cells.Item([int]$row, 4).Value = $item.D
    }
    if ($item.ContainsKey("E")) {
        $ws.Cells.Item([int]$row, 5).Value = $item.E
    }
}
